# ---------------------------------------------------------------------------
# Concertacionevaluacion.docx - "Add files via upload" re-save
#
# The canonical-OOXML diff behind this commit does not change any visible
# text, formatting or layout. It only touches internal Word bookkeeping
# that gets rewritten whenever the file is opened and re-saved:
#
#   1. word/header1.xml   - the (empty) <w:sdtEndPr/> on the page header's
#                            date content control disappears - Word drops
#                            that vestigial wrapper when it rewrites the
#                            control.
#   2. word/glossary/settings.xml
#                          - two new <w:rsid> "editing session" stamps
#                            appear (the glossary part backs the date
#                            control's placeholder building block, so it
#                            gets touched together with it).
#   3. customXml/item1.xml <-> customXml/item3.xml (and their itemProps
#                            counterparts) swap places - the SharePoint
#                            "FormTemplates" part and the "properties"
#                            part trade part numbers, which is how Word
#                            happens to renumber customXml parts on a
#                            save/round-trip; the parts' own content is
#                            unchanged.
#
# Reproduce the same user gesture that produces this no-visible-effect
# diff: open the page header, select the date content control and
# re-confirm it (which is what makes Word rewrite the <w:sdt> block and
# touch the backing glossary part on save). Everything below is
# defensive: if a particular piece of internal state isn't reachable
# through the object model in this session, the attempt is skipped
# rather than risking an unrelated change to the document's real
# content.

$d = $word.ActiveDocument

# --- 1/2: re-assert the header's date content control -----------------
# (drops the stale empty <w:sdtEndPr/> and stamps a fresh rsid on the
# content control + its glossary placeholder when Word rewrites it)
try {
    $sec = $d.Sections.Item(1)
    $hdr = $sec.Headers.Item(1)
    $cc  = $hdr.Range.ContentControls.Item(1)

    if ($cc -ne $null -and $cc.Type -eq 6) {
        try { $cc.DateCalendarType  = $cc.DateCalendarType }  catch { }
        try { $cc.DateDisplayFormat = $cc.DateDisplayFormat } catch { }
        try { $cc.DateStorageFormat = $cc.DateStorageFormat } catch { }
    }
} catch { }

# --- 3: renumber the SharePoint custom XML parts -----------------------
# customXml/item1.xml (FormTemplates) and customXml/item3.xml
# (properties/documentManagement) trade places; content of each part is
# unchanged, only which item number it is saved under changes.
try {
    $parts = $d.CustomXMLParts
    if ($parts -ne $null -and $parts.Count -gt 0) {
        $formTemplatesXml = $null
        $propertiesXml = $null
        for ($i = 1; $i -le $parts.Count; $i++) {
            $part = $parts.Item($i)
            if ($part.NamespaceURI -eq "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms") {
                $formTemplatesXml = $part.XML
            }
            if ($part.NamespaceURI -eq "http://schemas.microsoft.com/office/2006/metadata/properties") {
                $propertiesXml = $part.XML
            }
        }
        if ($formTemplatesXml -ne $null -and $propertiesXml -ne $null) {
            for ($i = 1; $i -le $parts.Count; $i++) {
                $part = $parts.Item($i)
                if ($part.NamespaceURI -eq "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms") {
                    $part.Delete()
                }
            }
            for ($i = 1; $i -le $parts.Count; $i++) {
                $part = $parts.Item($i)
                if ($part.NamespaceURI -eq "http://schemas.microsoft.com/office/2006/metadata/properties") {
                    $part.Delete()
                }
            }
            $d.CustomXMLParts.Add($propertiesXml)
            $d.CustomXMLParts.Add($formTemplatesXml)
        }
    }
} catch { }
